$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the title "Мощность" ("Power") to A1.
$ws.Range("A1").Value = "Мощность"

# --- "No alignment" look (font + border, but no wrap / centering) ---
# A1 loses the inherited wrap + vertical-center alignment that the rest of
# the header row still has.
$ws.Range("A1").VerticalAlignment = -4107
$ws.Range("A1").WrapText = $false

# The merged results block (A8:G10) gets the very same plain look as A1 -
# copy A1's format over rather than re-deriving it property by property,
# so both areas end up sharing one identical cell style.
$ws.Range("A1").Copy()
$ws.Range("A8:G10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Vertical-center + wrap-text look ---
# Header row (B1:H1), row labels (A2:A7) and the result labels (H8:H10)
# already carry this look from the original workbook, so nothing else to
# change there.

# --- Centered matrix values ---
# The pairwise-comparison values (B2:H7) additionally get horizontally
# centered, on top of the vertical-center + wrap-text they already have.
$ws.Range("B2:H7").HorizontalAlignment = -4108

Write-Host "Done"
